# Ajustando a geraçao do relatorio
#
# 1) "nome_voluntario" placeholder textbox: underline just the
#    "voluntario" part (splits the single run into "nome_" + "voluntario").
# 2) The certificate body textbox: drop the << >> markers around
#    NOME_ONG and append " Horas" after CARGA_HORARIA.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape "Google Shape;103;p13" -> nome_voluntario ---
$nameShape = $s.Shapes.Item(12)
$nameRange = $nameShape.TextFrame.TextRange
$underlinePart = $nameRange.Characters(6, 10)
$underlinePart.Font.Underline = $true

# --- Shape "Google Shape;104;p13" -> certificate body text ---
$bodyShape = $s.Shapes.Item(13)
$bodyShape.TextFrame.TextRange.Text = "PARTICIPOU DE UMA AÇÃO SOCIAL DA ONG  NOME_ONG, COM A CARGA HORÁRIA DE CARGA_HORARIA Horas."
